$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-text numbers (some using dotted
# thousands separators, trailing zeros, or subscript digits), so force
# text formatting on the price cells we touch before writing their new
# values -- this stops Excel from silently re-parsing them as floats.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($pc in $priceCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.208.49"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.505.74"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "199.33"
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("D6").Value = "547.18"
$ws.Range("E6").Value = "  -6.42%  "
$ws.Range("D7").Value = "3.490.28"
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("D8").Value = "0.602"
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.652"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "63.31"
$ws.Range("E11").Value = "  +14.89%  "
$ws.Range("D12").Value = "0.142"
$ws.Range("E12").Value = "  -5.93%  "
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  -7.13%  "
$ws.Range("D14").Value = "9.76"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "4.093.91"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "3.525.72"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "18.44"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "67.159.27"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  -6.18%  "
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("D22").Value = "389.10"
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("D23").Value = "3.99"
$ws.Range("E23").Value = "  -6.06%  "
$ws.Range("D24").Value = "11.87"
$ws.Range("E24").Value = "  -11.47%  "
$ws.Range("D25").Value = "82.17"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("D26").Value = "12.15"
$ws.Range("E26").Value = "  -3.20%  "
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  -5.13%  "
$ws.Range("D28").Value = "3.73"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").Value = "8.79"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("D30").Value = "30.75"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("D31").Value = "684.84"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "7.14"
$ws.Range("E32").Value = "  -11.68%  "
$ws.Range("D33").Value = "11.72"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").Value = "63.63"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("D36").Value = "38.55"
$ws.Range("E36").Value = "  -9.42%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "0.402"
$ws.Range("E38").Value = "  -4.94%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.070.23"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "3.00"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("D43").Value = "0.0₃0678"
$ws.Range("E43").Value = "  -13.38%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  -15.03%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  +8.85%  "
$ws.Range("D46").Value = "2.72"
$ws.Range("E46").Value = "  +6.27%  "
$ws.Range("D47").Value = "0.0397"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").Value = "138.68"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").Value = "8.24"
$ws.Range("E50").Value = "  -6.03%  "
$ws.Range("D51").Value = "2.90"
$ws.Range("E51").Value = "  -6.17%  "
